# The worksheet starts as a small table:
#   Category | Overlap | Total | Percentage
#   police   | 0       | 3     | 0
#   schools  | 0       | 10    | 0
#
# The edit inserts a new "health" row above "police" (pushing the existing
# rows down by one), gives "police" a fresh Total/Percentage of 3/0, while
# "schools" keeps its original Total/Percentage of 10/0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 4 first. It inherits formatting from row 3
# (an unstyled data row), so no extra style gets created, and it makes room
# to shift the existing data down without ever inserting next to the
# styled header row.
$ws.Rows.Item(4).Insert()

# Row 4 becomes the "schools" row (the original row 3 data).
$ws.Cells.Item(4, 1).Value = "schools"
$ws.Cells.Item(4, 2).Value = 0
$ws.Cells.Item(4, 3).Value = 10
$ws.Cells.Item(4, 4).Value = 0

# Row 3 becomes the "police" row, with its new Total/Percentage values.
$ws.Cells.Item(3, 1).Value = "police"
$ws.Cells.Item(3, 2).Value = 0
$ws.Cells.Item(3, 3).Value = 3
$ws.Cells.Item(3, 4).Value = 0

# Row 2 becomes the newly inserted "health" row.
$ws.Cells.Item(2, 1).Value = "health"
$ws.Cells.Item(2, 2).Value = 0
$ws.Cells.Item(2, 3).Value = 3
$ws.Cells.Item(2, 4).Value = 0
